$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 1.7
$ws.Range("I2").Value = 5.5
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 7
$ws.Range("AE2").Value = 21
$ws.Range("AN2").Value = 3.5

# Row 12 updates
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 3.25
